# Bump the "Förändrad" (Changed) date in column C by one day for every
# data row in the sheet (rows 2 through the last used row).
#
# Column C holds a date serial number (formatted as YYYY-MM-DD). In the
# source workbook every data row had the same value, 45188 (2023-09-19),
# which the commit updates uniformly to 45189 (2023-09-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
